$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, centered, bordered) from B1 into H1, then set its value
$ws.Range("B1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Save values per row (2..11)
$saveValues = @{
    2 = 0
    3 = 0
    4 = 1
    5 = 1
    6 = 0
    7 = 1
    8 = 1
    9 = 1
    10 = 1
    11 = 1
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
